# "Revert link to survey deleted" -- restores the previous input figures on
# the "intent_01" sheet. Only D1 (Cjm[NO2]), D14 (Fmeteo) and D18 (Road
# length) are independent inputs; every other changed cell (D3, D9, D10,
# D11, D15, D16, D17, D19) is a formula that recalculates automatically
# from these.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intent_01")

$ws.Range("D1").Value  = 15      # Cjm[NO2]    (was 5.13)
$ws.Range("D14").Value = 0.5     # Fmeteo      (was 1.8)
$ws.Range("D18").Value = 18075   # Road length (was 8003)

# Restore the saved view state for the sheet (zoom level + selected cell).
$ws.Activate()
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("C18").Select() | Out-Null
